# Update crypto price/volume table to the latest scraped snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.067.27"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "2.386.71"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'557.58"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").Value = "'134.14"
$ws.Range("E6").Value = "  -1.95%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.65%  "

$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D12").Value = "'0.344"
$ws.Range("E12").Value = "  -2.60%  "

$ws.Range("D13").Value = "'24.51"
$ws.Range("E13").Value = "  -3.03%  "

$ws.Range("D14").Value = "2.811.28"
$ws.Range("E14").Value = "  -0.79%  "

$ws.Range("D15").Value = "60.017.96"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("D17").Value = "2.390.24"
$ws.Range("E17").Value = "  -0.65%  "

$ws.Range("D18").Value = "'11.12"
$ws.Range("E18").Value = "  -1.57%  "

$ws.Range("D19").Value = "'4.50"
$ws.Range("E19").Value = "  +2.15%  "

$ws.Range("D20").Value = "'321.75"
$ws.Range("E20").Value = "  -1.98%  "

$ws.Range("D21").Value = "'6.68"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "'64.12"
$ws.Range("E23").Value = "  -3.57%  "

$ws.Range("D24").Value = "'0.174"
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").Value = "'8.46"
$ws.Range("E26").Value = "  -2.26%  "

$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("E28").Value = "  +2.46%  "

$ws.Range("D29").Value = "0.0₃0762"
$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("D30").Value = "'169.89"
$ws.Range("E30").Value = "  +0.92%  "

$ws.Range("D31").Value = "'6.10"
$ws.Range("E31").Value = "  +1.13%  "

$ws.Range("D32").Value = "'1.12"
$ws.Range("E32").Value = "  +11.06%  "

$ws.Range("D33").Value = "'0.401"
$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("D34").Value = "'18.18"
$ws.Range("E34").Value = "  -2.26%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.33"
$ws.Range("E35").Value = "  +2.21%  "

$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").Value = "'4.16"
$ws.Range("E38").Value = "  -0.54%  "

$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").Value = "'319.93"
$ws.Range("E40").Value = "  +0.41%  "

$ws.Range("D41").Value = "'38.67"
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("D42").Value = "'145.54"
$ws.Range("E42").Value = "  +4.32%  "

$ws.Range("E43").Value = "  -3.12%  "

$ws.Range("D44").Value = "'0.0971"
$ws.Range("E44").Value = "  +0.51%  "

$ws.Range("D45").Value = "'19.83"
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("D46").Value = "'0.0513"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").Value = "'0.0219"
$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("D49").Value = "'11.06"
$ws.Range("E49").Value = "  +0.28%  "

$ws.Range("D50").Value = "'1.55"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("D51").Value = "'4.67"
$ws.Range("E51").Value = "  -0.08%  "

